$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 264, shifting existing rows 264-337 down to 265-338.
$ws.Rows.Item(264).Insert()

# Populate the newly inserted row 264 with its data. Most fields mirror the
# row that used to occupy this position (now at row 265): A, B, C, E, F, G,
# H, I, J, N, O, Q, R stay the same; D, K, L, M, P get the new values.
$ws.Cells.Item(264, 1).Value = 6
$ws.Cells.Item(264, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(264, 3).Value = "Metropolitana"
$ws.Cells.Item(264, 4).Value = 45093
$ws.Cells.Item(264, 5).Value = 13
$ws.Cells.Item(264, 6).Value = 100112029
$ws.Cells.Item(264, 7).Value = "Orégano"
$ws.Cells.Item(264, 8).Value = "Sin especificar"
$ws.Cells.Item(264, 9).Value = "Primera"
$ws.Cells.Item(264, 10).Value = 35
$ws.Cells.Item(264, 11).Value = 19000
$ws.Cells.Item(264, 12).Value = 20000
$ws.Cells.Item(264, 13).Value = 19571
$ws.Cells.Item(264, 14).Value = "$/docena de atados"
$ws.Cells.Item(264, 15).Value = "Región Metropolitana"
$ws.Cells.Item(264, 16).Value = 6524
$ws.Cells.Item(264, 17).Value = 3
$ws.Cells.Item(264, 18).Value = "Hortaliza"
